# Commit #5: insurance, claim, debt, investment done
# This adds proper headers + trailing metadata columns (property_category,
# category, date, legislator_name, legislator_id, source_file, index) to the
# "債務" (debt) worksheet (sheet6), matching the convention already used by
# the other sheets in this workbook, and relabels the mis-used header row
# (which previously held a copy of row 2's data) with real column names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)

# --- Row 1: header labels -------------------------------------------------
$ws.Range("B1").Value = "species"
$ws.Range("C1").Value = "debtor"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Give the new header cells (H1:N1) the same look as the existing header
# cells (bold, centered/top-aligned, thin box border) used for B1:G1.
$headerNew = $ws.Range("H1:N1")
$headerNew.Font.Bold = $true
$headerNew.HorizontalAlignment = -4108
$headerNew.VerticalAlignment = -4160
$headerNew.Borders.LineStyle = 1

# --- Row 2 (record #120) ---------------------------------------------------
$ws.Range("H2").Value = "debt"
$ws.Range("I2").Value = "normal"
# Force the "date" column to be stored as plain text (it mirrors the
# "2012-04-27" source-file date, not a real Excel date value) instead of
# letting Excel auto-convert it to a date serial number.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2012-04-27"
$ws.Range("J2").ClearFormats()
$ws.Range("K2").Value = "江惠貞"
$ws.Range("L2").Value = 1732
$ws.Range("M2").Value = "tmpf6f41"
$ws.Range("N2").Value = 120

# --- Row 3 (record #121) ---------------------------------------------------
$ws.Range("H3").Value = "debt"
$ws.Range("I3").Value = "normal"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "2012-04-27"
$ws.Range("J3").ClearFormats()
$ws.Range("K3").Value = "江惠貞"
$ws.Range("L3").Value = 1732
$ws.Range("M3").Value = "tmpf6f41"
$ws.Range("N3").Value = 121
